$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 102
$ws1.Range("F3").Value = 306

# Same update mirrored on the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 102
$ws4.Range("F3").Value = 306
